# Refresh the cryptocurrency price/volume snapshot (Price = col D, Volume(1h) = col E).
# Values are entered with a leading apostrophe so Excel stores them as literal text
# (matching the sheet's existing inline-string cells) instead of re-interpreting them
# as numbers/percentages, which would silently drop trailing zeros (e.g. "0.1290" -> 0.129)
# or flip tiny values into scientific notation (e.g. "0.00006984" -> 6.984E-05).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.77"
$ws.Range("E2").Value = "'-0.41%"
$ws.Range("D3").Value = "'41.65"
$ws.Range("E3").Value = "'1.08%"
$ws.Range("D4").Value = "'5.686"
$ws.Range("E4").Value = "'-0.55%"
$ws.Range("D5").Value = "'0.08345"
$ws.Range("E5").Value = "'3.32%"
$ws.Range("D6").Value = "'8.807"
$ws.Range("E6").Value = "'0.83%"
$ws.Range("D7").Value = "'1.993"
$ws.Range("E7").Value = "'-1.97%"
$ws.Range("D8").Value = "'4.468"
$ws.Range("E8").Value = "'-1.55%"
$ws.Range("D10").Value = "'0.9261"
$ws.Range("E10").Value = "'0.38%"
$ws.Range("D11").Value = "'0.1290"
$ws.Range("E11").Value = "'2.12%"
$ws.Range("D12").Value = "'0.1972"
$ws.Range("E12").Value = "'0.93%"
$ws.Range("D13").Value = "'0.09471"
$ws.Range("E13").Value = "'1.30%"
$ws.Range("D14").Value = "'0.03853"
$ws.Range("E14").Value = "'4.31%"
$ws.Range("D15").Value = "'0.1060"
$ws.Range("E15").Value = "'0.72%"
$ws.Range("D16").Value = "'0.001308"
$ws.Range("E16").Value = "'0.14%"
$ws.Range("D17").Value = "'0.006110"
$ws.Range("E17").Value = "'-2.83%"
$ws.Range("E18").Value = "'1.86%"
$ws.Range("E19").Value = "'1.54%"
$ws.Range("D20").Value = "'8.693"
$ws.Range("E20").Value = "'-1.92%"
$ws.Range("E21").Value = "'-3.98%"
$ws.Range("D22").Value = "'0.2486"
$ws.Range("E22").Value = "'-6.46%"
$ws.Range("D23").Value = "'0.04417"
$ws.Range("E23").Value = "'-0.13%"
$ws.Range("D24").Value = "'0.001279"
$ws.Range("E24").Value = "'1.32%"
$ws.Range("D25").Value = "'0.004382"
$ws.Range("E25").Value = "'1.69%"
$ws.Range("D26").Value = "'0.0001221"
$ws.Range("D39").Value = "'0.02837"
$ws.Range("E39").Value = "'-1.23%"
$ws.Range("D40").Value = "'0.05512"
$ws.Range("E40").Value = "'0.48%"
$ws.Range("D41").Value = "'0.007956"
$ws.Range("E41").Value = "'2.36%"
$ws.Range("D42").Value = "'0.1433"
$ws.Range("E42").Value = "'1.02%"
$ws.Range("D43").Value = "'0.009300"
$ws.Range("E43").Value = "'-6.36%"
$ws.Range("D44").Value = "'0.002142"
$ws.Range("E44").Value = "'-4.29%"
$ws.Range("D45").Value = "'0.01170"
$ws.Range("E45").Value = "'5.72%"
$ws.Range("D46").Value = "'0.00006984"
$ws.Range("E46").Value = "'2.18%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.24%"
$ws.Range("E48").Value = "'14.14%"
$ws.Range("D49").Value = "'0.002280"
$ws.Range("E49").Value = "'-0.30%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.24%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.24%"
